# Update the Gantt Chart worksheet with the latest project tracking data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 16: Champ. Player Special Abilities Dev -> Player Drop Bomb Ability ---
$ws.Range('C16').Value = 'Player Drop Bomb Ability'
$ws.Range('E16').Value = 45737
$ws.Range('F16').Value = 45744
$ws.Range('H16').Value = 'Complete'
$ws.Range('I16').Value = 1

# --- Row 17: Bomb animation -> now Complete ---
$ws.Range('H17').Value = 'Complete'
$ws.Range('I17').Value = 1

# --- Row 20: Bomb destruction -> now Complete ---
$ws.Range('H20').Value = 'Complete'
$ws.Range('I20').Value = 1

# --- Row 21: Merging and Status Update Meeting, date slip + completed ---
$ws.Range('E21').Value = 45743
$ws.Range('F21').Value = 45743
$ws.Range('H21').Value = 'Complete'
$ws.Range('I21').Value = 1
# Clear the one-off cyan highlight previously on C21 so it matches the other rows
$ws.Range('C21').Interior.Pattern = -4142

# --- Row 22: Comment Background and Tile Manager Class -> now Complete ---
$ws.Range('H22').Value = 'Complete'
$ws.Range('I22').Value = 1

# --- Row 23: Special door Dev -> Special Door Dev, reschedule + in progress ---
$ws.Range('C23').Value = 'Special Door Dev'
$ws.Range('E23').Value = 45742
$ws.Range('F23').Value = 45748
$ws.Range('H23').Value = 'In progress'
$ws.Range('I23').Value = 0.6

# --- Row 24: Comment all Enemy Classes -> now In progress ---
$ws.Range('H24').Value = 'In progress'
$ws.Range('I24').Value = 0.1

# --- Row 25: Comment Character and JackBomber Classes -> now In progress ---
$ws.Range('H25').Value = 'In progress'
$ws.Range('I25').Value = 0.15

# --- Row 26: new activity "Player Health Bar " ---
$ws.Range('C26').Value = 'Player Health Bar '
$ws.Range('E26').Value = 45744
$ws.Range('F26').Value = 45751
$ws.Range('H26').Value = 'In progress'
$ws.Range('I26').Value = 0.01

# --- Row 27: new activity "Game Level Timer and Bar" ---
$ws.Range('C27').Value = 'Game Level Timer and Bar'
$ws.Range('D27').Value = 'Murat C. GZ'
$ws.Range('E27').Value = 45744
$ws.Range('F27').Value = 45751

# --- Row 28: new activity "Game Entry GUI" ---
$ws.Range('C28').Value = 'Game Entry GUI'
$ws.Range('D28').Value = 'Andrias'
$ws.Range('E28').Value = 45744
$ws.Range('F28').Value = 45751

# --- Row 29: new activity "Player  Bomb Drop Limitiatiions by Timer" ---
$ws.Range('C29').Value = 'Player  Bomb Drop Limitiatiions by Timer'
$ws.Range('D29').Value = 'Murat C. GZ'
$ws.Range('E29').Value = 45751
$ws.Range('F29').Value = 45758

# --- Row 30: new activity "Game Over Animations" ---
$ws.Range('C30').Value = 'Game Over Animations'
$ws.Range('D30').Value = 'Andrias'
$ws.Range('E30').Value = 45751
$ws.Range('F30').Value = 45758
$ws.Range('H30').Value = 'Not Started'
$ws.Range('I30').Value = 0
$ws.Range('G30').Formula = '=IF(F30="","",NETWORKDAYS(E30,F30))'

# --- Row 31: new schedule slot (no activity name yet) ---
$ws.Range('E31').Value = 45753
$ws.Range('F31').Value = 45760
$ws.Range('H31').Value = 'Not Started'
$ws.Range('I31').Value = 0
$ws.Range('G31').Formula = '=IF(F31="","",NETWORKDAYS(E31,F31))'

# --- Row 32: new schedule slot ---
$ws.Range('E32').Value = 45754
$ws.Range('F32').Value = 45761
$ws.Range('H32').Value = 'Not Started'
$ws.Range('I32').Value = 0
$ws.Range('G32').Formula = '=IF(F32="","",NETWORKDAYS(E32,F32))'

# --- Row 33: new schedule slot ---
$ws.Range('E33').Value = 45755
$ws.Range('F33').Value = 45762
$ws.Range('H33').Value = 'Not Started'
$ws.Range('I33').Value = 0
$ws.Range('G33').Formula = '=IF(F33="","",NETWORKDAYS(E33,F33))'

# --- Row 34: new schedule slot ---
$ws.Range('E34').Value = 45756
$ws.Range('F34').Value = 45763
$ws.Range('H34').Value = 'Not Started'
$ws.Range('I34').Value = 0
$ws.Range('G34').Formula = '=IF(F34="","",NETWORKDAYS(E34,F34))'

# --- Row 35: new schedule slot ---
$ws.Range('E35').Value = 45757
$ws.Range('F35').Value = 45764
$ws.Range('H35').Value = 'Not Started'
$ws.Range('I35').Value = 0
$ws.Range('G35').Formula = '=IF(F35="","",NETWORKDAYS(E35,F35))'

# --- Row 36: new schedule slot ---
$ws.Range('E36').Value = 45758
$ws.Range('F36').Value = 45765
$ws.Range('H36').Value = 'Not Started'
$ws.Range('I36').Value = 0
$ws.Range('G36').Formula = '=IF(F36="","",NETWORKDAYS(E36,F36))'

# --- Row 37: new schedule slot ---
$ws.Range('E37').Value = 45759
$ws.Range('F37').Value = 45766
$ws.Range('H37').Value = 'Not Started'
$ws.Range('I37').Value = 0
$ws.Range('G37').Formula = '=IF(F37="","",NETWORKDAYS(E37,F37))'

# --- Row 38: new schedule slot ---
$ws.Range('E38').Value = 45760
$ws.Range('F38').Value = 45767
$ws.Range('H38').Value = 'Not Started'
$ws.Range('I38').Value = 0
$ws.Range('G38').Formula = '=IF(F38="","",NETWORKDAYS(E38,F38))'

# --- Final selection matches the saved view in the source workbook ---
$ws.Range('F32').Select()
